$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("C4").Value = "RW-02, PC-01"
$ws.Range("C4").Select()
